# Fixed test template version: remove the stale "M2Doc version mismatch"
# warning runs that were left in the template between the spell-check
# markers, right before the "query" word.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    query",
    $true, $false, $false, $false, $false, $true, 1, $false, "query", 2)
